$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.799.81"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "1.574.23"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.63%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("E6").Value = "  +6.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "25.37"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.65%  "
$ws.Range("E9").Value = "  +2.25%  "
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.33%  "
$ws.Range("D12").Value = "1.798.13"
$ws.Range("E12").Value = "  +1.46%  "
$ws.Range("D13").Value = "1.577.22"
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").Value = "28.796.96"
$ws.Range("E15").Value = "  +2.05%  "
$ws.Range("E16").Value = "  +1.34%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "230.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.57%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("E20").Value = "  +2.74%  "
$ws.Range("E21").Value = "  -0.78%  "
$ws.Range("E22").Value = "  +1.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.26"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.108"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.89"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.25%  "
$ws.Range("E28").Value = "  +1.64%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0461"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.79%  "
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("E32").Value = "  +1.14%  "
$ws.Range("D33").Value = "1.415.75"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.04"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.76"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.17%  "
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0162"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.07%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.521"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.52%  "
$ws.Range("E41").Value = "  +1.76%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.781"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("E44").Value = "  -0.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "64.12"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.17%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").Value = "1.710.16"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.832"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.53%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "84.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.46"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.54%  "
